$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vars_meta_data")
# Sanity: this is the sheet that was active/edited in the source workbook
# (equivalent to $wb.ActiveSheet for this file).

# Row 41 - Longhurst_Long
$ws.Range("A41").Value = "Longhurst_Long"
$ws.Range("B41").Value = "Longhurst province sample was taken in."
$ws.Range("C41").Value = "NA"
$ws.Range("D41").Value = "NA"
$ws.Range("E41").Value = "Irregular"
$ws.Range("F41").Value = "Irregular"
$ws.Range("G41").Value = "Biology"
$ws.Range("H41").Value = 1

# Row 42 - Longhurst_Short
$ws.Range("A42").Value = "Longhurst_Short"
$ws.Range("B42").Value = "Longhurst province sample was taken in, shortened code."
$ws.Range("C42").Value = "NA"
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = "Irregular"
$ws.Range("F42").Value = "Irregular"
$ws.Range("G42").Value = "Biology"
$ws.Range("H42").Value = 1

# Row 43 - Season
$ws.Range("A43").Value = "Season"
$ws.Range("B43").Value = "Season sample was taken in."
$ws.Range("C43").Value = "NA"
$ws.Range("D43").Value = "NA"
$ws.Range("E43").Value = "Irregular"
$ws.Range("F43").Value = "Irregular"
$ws.Range("G43").Value = "Biology"
$ws.Range("H43").Value = 1

# Update selection / active cell to match the edited range
$ws.Activate() | Out-Null
$ws.Range("A41:H43").Select() | Out-Null
